$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New unique-values were appended into several of the stacked "array" blocks
# in column B (Price, MPG, Color, Horsepower). Each block is visually grouped
# by a fill style; inserting a row within a block (shift down, format copied
# from the row above) keeps that block's styling intact.
#
# Work from the bottom of the sheet upward so row numbers for not-yet-processed
# insertions stay valid (everything below the insertion point simply shifts).
# ---------------------------------------------------------------------------

# Horsepower block: <300 (row22), <500 (row23) -> <200, <300, <400, <500
$ws.Rows("23:23").Insert() | Out-Null
$ws.Range("B23").Value = "<400"

$ws.Rows("22:22").Insert() | Out-Null
$ws.Range("B22").Value = "<200"

# Color block: red(16) black(17) silver(18) -> add "blue" after silver
$ws.Rows("19:19").Insert() | Out-Null
$ws.Range("B19").Value = "blue"

# MPG block: <20(8) <30(9) <40(10) -> add <25 after <20, <35 after <30
$ws.Rows("10:10").Insert() | Out-Null
$ws.Range("B10").Value = "<35"

$ws.Rows("9:9").Insert() | Out-Null
$ws.Range("B9").Value = "<25"

# Price block: <20,000(3) <30,000(4) <40,000(5) <50,000(6) <60,000(7)
# -> add <25,000 / <35,000 / <45,000 / <55,000
$ws.Rows("7:7").Insert() | Out-Null
$ws.Range("B7").Value = "<55,000"

$ws.Rows("6:6").Insert() | Out-Null
$ws.Range("B6").Value = "<45,000"

$ws.Rows("5:5").Insert() | Out-Null
$ws.Range("B5").Value = "<35,000"

$ws.Rows("4:4").Insert() | Out-Null
$ws.Range("B4").Value = "<25,000"

# The new Horsepower row ("<200", now row 29) inherited the fill style of the
# row above it (the "Engine Fuel Type" block) because it was inserted right
# after "diesel". Fix its fill to match the rest of the Horsepower block by
# copying the format from "<300" (row 30), which already carries the correct
# style.
$ws.Range("B30").Copy() | Out-Null
$ws.Range("B29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Renumber the sequential "index" column (A3:A44 = 0..41) now that new rows
# have been woven into the list.
# ---------------------------------------------------------------------------
for ($i = 0; $i -le 41; $i++) {
    $ws.Cells.Item(3 + $i, 1).Value = $i
}

# Match the saved selection/active cell from the authored workbook.
$ws.Range("A45").Select() | Out-Null
